# curso_elton.xlsx - "Add files via upload" edit
#
# 1. Plan1!F64 comment gets two extra lines inserted right after the
#    "Qual o melhor formato..." question.
# 2. Plan1 rows 63 and 65 gain Assistido/Implementado (sim/não) marks plus a
#    motivational "Mensagem" quote in column G (two brand-new shared
#    strings).
# 3. The active selection on Plan1 moves from F26 (A26:F26) down to A62,
#    with the view scrolled so row 41 is back at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- 1. Extend the F64 comment left by Alan Jose do Nascimento ---------
$comment = $ws.Range("F64").Comment
$newCommentText = @"
Alan Jose do Nascimento:
Seguem os links das dicas:

VideoMaker da Ozi: http://eu.queroumvideomaker.com/

Meu Editor: Heider Freitas 61 98507 1172

Qual o melhor formato de experiementação que eu gosto e posso gerar?

Qual tipo de formato que o meu publico gosta?
A experimentação deve ser no nosso ambiente, ou ambiente reservado.

Separar as 3 principais ideias do papo que cola colocar no espaco de ideias do canvas 

Experimentacao tem que acontecer no meu ambiente;

Definir o formato da entrega das 3 ideias (vides, ebook, palestra etc)

"@
$comment.Text($newCommentText)

# --- 2. Fill in the "Assistido" / "Implementado" / "Mensagem" columns --
# Row 63 ("A4. Pensando em formas de atrair as pessoas")
$ws.Range("C63").Value = "sim"
$ws.Range("D63").Value = "não"
$ws.Range("G63").Value = "Usar as mesmas palavras não é garantia de entendimento. É preciso ter experiência em comum com alguém. Friedrich Nietzche"

# Row 65 ("A6. Afinando o fechamento")
$ws.Range("C65").Value = "sim"
$ws.Range("D65").Value = "não"
$ws.Range("G65").Value = "Para vencer na vida você não precisa fazer o que os outros não fazem, precisa fazer o que você não faz."

# --- 3. Move the selection / scroll position ----------------------------
$ws.Activate()
$ws.Range("A62").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 41
$window.ScrollColumn = 1
